$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure text-valued cells (prices/links/names) are not re-interpreted as numbers

# Row 2
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "90.958.99"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "  +4.27%  "

# Row 3
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "3.091.26"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "  +1.22%  "

# Row 4
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "  +0.10%  "

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "218.78"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "  +4.51%  "

# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "618.77"
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "  -0.52%  "

# Row 7
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.377"
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "  +4.31%  "

# Row 8
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "  +15.13%  "

# Row 9
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "  +0.03%  "

# Row 10
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "3.088.43"
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "  +1.20%  "

# Row 11
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.700"
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "  +20.47%  "

# Row 12
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.190"
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "  +7.95%  "

# Row 13
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "  +8.70%  "

# Row 14
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "5.38"
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "  +2.88%  "

# Row 15
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "90.722.12"
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "  +4.34%  "

# Row 16
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "33.05"
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "  +5.91%  "

# Row 17
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "3.669.27"
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "  +1.70%  "

# Row 18
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "3.085.79"
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "  +0.37%  "

# Row 19
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "3.62"
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = "  +7.77%  "

# Row 20
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "0.0000234"
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "  +13.18%  "

# Row 21
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "  +6.56%  "

# Row 22
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "431.20"
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "  +4.07%  "

# Row 23
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "8.56"
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "  +4.86%  "

# Row 24
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "5.13"
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = "  +7.64%  "

# Row 25
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "5.56"
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "  +2.32%  "

# Row 26
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "11.88"
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = "  +5.80%  "

# Row 27
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "83.66"
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = "  +2.17%  "

# Row 28
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "3.262.20"
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = "  +1.39%  "

# Row 29
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "0.997"
$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = "  -0.24%  "

# Row 30
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "0.167"
$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = "  +12.44%  "

# Row 31
$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = "  +0.04%  "

# Row 32
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "8.68"
$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = "  +8.79%  "

# Row 33
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "3.90"
$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = "  +8.38%  "

# Row 34
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "517.80"
$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = "  +4.55%  "

# Row 35
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "6.93"
$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = "  +5.31%  "

# Row 36
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.139"
$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = "  -0.76%  "

# Row 37
$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = "  +4.46%  "

# Row 38
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = "  +3.18%  "

# Row 39
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "22.99"
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = "  +5.51%  "

# Row 40
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "22.29"
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "  +0.71%  "

# Row 41
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "  +0.12%  "

# Row 42
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "  +11.72%  "

# Row 43
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "  +0.01%  "

# Row 44
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.371"
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "  +3.38%  "

# Row 45
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = "  +3.98%  "

# Row 46
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.0721"
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = "  +13.04%  "

# Row 47
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "43.88"
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = "  +1.16%  "

# Row 48
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "141.54"
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = "  -3.45%  "

# Row 49
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = "  +8.52%  "

# Row 50
$ws.Range("B50").Value = "Filecoin"
$ws.Range("C50").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "4.19"
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = "  +8.64%  "

# Row 51
$ws.Range("B51").Value = "FLOKI"
$ws.Range("C51").Value = "https://coinranking.com/coin/fmHk13Rqw+floki-floki"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.000258"
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = "  +19.22%  "
